$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.200.87'
$ws.Range("E2").Value = '  +8.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.588.79'
$ws.Range("E3").Value = '  +7.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9906'
$ws.Range("E5").Value = '  +2.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.08'
$ws.Range("E6").Value = '  +7.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3615'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  +8.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.29'
$ws.Range("E9").Value = '  +4.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.115'
$ws.Range("E10").Value = '  +3.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06932'
$ws.Range("E11").Value = '  +4.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.38'
$ws.Range("E13").Value = '  +6.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.814'
$ws.Range("E14").Value = '  +5.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.515'
$ws.Range("E15").Value = '  +5.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9908'
$ws.Range("E16").Value = '  +3.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001061'
$ws.Range("E17").Value = '  +3.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.585.29'
$ws.Range("E18").Value = '  +7.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06588'
$ws.Range("E19").Value = '  +10.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.90'
$ws.Range("E20").Value = '  +10.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.78'
$ws.Range("E21").Value = '  +8.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.900'
$ws.Range("E22").Value = '  +7.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.57'
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.187.86'
$ws.Range("E24").Value = '  +7.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.371'
$ws.Range("E25").Value = '  +4.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.489'
$ws.Range("E26").Value = '  +17.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.07'
$ws.Range("E27").Value = '  +3.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.10'
$ws.Range("E28").Value = '  +11.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.755.90'
$ws.Range("E29").Value = '  +7.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.85'
$ws.Range("E30").Value = '  +5.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.917'
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.863'
$ws.Range("E32").Value = '  +18.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9166'
$ws.Range("E33").Value = '  +14.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08132'
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.626'
$ws.Range("E35").Value = '  +7.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.67'
$ws.Range("E36").Value = '  +11.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.107'
$ws.Range("E37").Value = '  +7.93%  '
$ws.Range("E38").Value = '  +1.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06021'
$ws.Range("E39").Value = '  +4.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.301'
$ws.Range("E40").Value = '  +11.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02177'
$ws.Range("E41").Value = '  +5.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1975'
$ws.Range("E42").Value = '  +5.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9905'
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5757'
$ws.Range("E44").Value = '  +9.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.762'
$ws.Range("E45").Value = '  +6.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.85'
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '124.70'
$ws.Range("E47").Value = '  +4.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5549'
$ws.Range("E48").Value = '  +6.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.929'
$ws.Range("E49").Value = '  +5.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06703'
$ws.Range("E50").Value = '  +3.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.08'
$ws.Range("E51").Value = '  +7.32%  '
